$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = [double]"0.06002781679967685"
$ws.Range("E2").Value = [double]"0.06002781679967685"

# Row 3
$ws.Range("D3").Value = [double]"1.368178241691181E-22"
$ws.Range("E3").Value = [double]"1.368178241691181E-22"

# Row 4
$ws.Range("D4").Value = [double]"0.001613110215634427"
$ws.Range("E4").Value = [double]"0.001613110215634427"

# Row 5
$ws.Range("D5").Value = [double]"0.8510013237162374"
$ws.Range("E5").Value = [double]"0.8510013237162374"

# Row 6
$ws.Range("D6").Value = [double]"0.9999430540994376"
$ws.Range("E6").Value = [double]"0.9999430540994376"

# Row 7
$ws.Range("D7").Value = [double]"0.9999999833317851"
$ws.Range("E7").Value = [double]"1.666821491941306E-08"

# Row 8
$ws.Range("D8").Value = [double]"0.9750024959055397"
$ws.Range("E8").Value = [double]"0.02499750409446033"

# Row 9
$ws.Range("D9").Value = [double]"0.9958556929827005"
$ws.Range("E9").Value = [double]"0.004144307017299509"

# Row 10
$ws.Range("D10").Value = [double]"1"
$ws.Range("E10").Value = [double]"0"

# Row 11
$ws.Range("D11").Value = [double]"0.9996056178377913"
$ws.Range("E11").Value = [double]"0.0003943821622086618"
$ws.Range("F11").Value = [double]"1.177060961723328"
